# Add a new row (19) to Sheet1 with a new DataCamp course entry, matching
# the formatting already used for the last few rows of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- New data row -----------------------------------------------------
$ws.Range("A19").Value = "Writing Functions in Python"
$ws.Range("B19").Value = 3

# --- Match the formatting used on the row above (A18 / B18) -----------
# Copy A18's format (colored header-like text style) onto A19.
$ws.Range("A18").Copy() | Out-Null
$ws.Range("A19").PasteSpecial(-4122) | Out-Null

# Match B18's font color/style onto B19 (automatic/theme text color, same
# as the number column uses everywhere else in the table).
$ws.Range("B19").Font.ThemeColor = 1

$excel.CutCopyMode = 0

# --- Keep selection / active cell consistent with a freshly-typed row -
$ws.Range("A20").Select() | Out-Null
